# Commit message (informal Hindi/Hinglish): "null check daale-initialQuizeData,
# dataa hataayaa dono excel se" -> "added a null check for initialQuizData,
# removed the data from both excel [sheets]".
#
# The diff shows Sheet1 previously held a "user"/"result" header with four
# data rows (large JSON blobs in column B). The edit clears that data,
# leaving only a header row, and retitles the headers to "index"/"json".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear out all existing data on the sheet (rows 2-5 had the JSON payloads).
$ws.Cells.Clear()

# Rewrite the header row with the new column names.
$ws.Range("A1").Value = "index"
$ws.Range("B1").Value = "json"

$ws.Select()
$ws.Range("A1").Select()
